$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.178.94"
$ws.Range("E2").Value = "  -1.56%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.277.03"
$ws.Range("E3").Value = "  -1.75%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.80"
$ws.Range("E5").Value = "  +0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "179.09"
$ws.Range("E6").Value = "  -2.51%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  +4.57%  "

$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("E9").Value = "  -2.19%  "

$ws.Range("E10").Value = "  +1.18%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.403"
$ws.Range("E11").Value = "  -0.83%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.851.42"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("E13").Value = "  -3.50%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "66.216.22"
$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("E15").Value = "  -3.38%  "

$ws.Range("E16").Value = "  -1.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.280.75"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "435.51"
$ws.Range("E18").Value = "  -1.98%  "

$ws.Range("E19").Value = "  -2.25%  "

$ws.Range("E20").Value = "  -2.95%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.40"
$ws.Range("E21").Value = "  -4.01%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.71"
$ws.Range("E22").Value = "  -3.10%  "

$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.423.13"
$ws.Range("E24").Value = "  -1.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.507"
$ws.Range("E25").Value = "  -0.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.197"
$ws.Range("E26").Value = "  +2.13%  "

$ws.Range("E27").Value = "  -5.61%  "

$ws.Range("E28").Value = "  -1.90%  "

$ws.Range("E29").Value = "  +0.03%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.95"
$ws.Range("E30").Value = "  -1.34%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "22.29"
$ws.Range("E31").Value = "  -2.55%  "

$ws.Range("E32").Value = "  +0.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").Value = "  -2.22%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.59"
$ws.Range("E34").Value = "  -2.46%  "

$ws.Range("E35").Value = "  -3.80%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "158.65"
$ws.Range("E36").Value = "  -1.82%  "

$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "26.53"
$ws.Range("E38").Value = "  -3.46%  "

$ws.Range("E39").Value = "  -2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.775.11"
$ws.Range("E40").Value = "  -1.99%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.774"
$ws.Range("E41").Value = "  -2.18%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.32"
$ws.Range("E42").Value = "  -3.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.26"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("E44").Value = "  -3.24%  "

$ws.Range("E45").Value = "  -1.69%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "321.72"
$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.28"
$ws.Range("E47").Value = "  -2.71%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.23"
$ws.Range("E48").Value = "  -4.95%  "

$ws.Range("E49").Value = "  -2.13%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.103"
$ws.Range("E50").Value = "  +3.13%  "

$ws.Range("E51").Value = "  +0.02%  "
